$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header H1 from pipe_length to length_pipe
$ws.Range("H1").Value = "length_pipe"

# Update numeric values in row 2 to reflect updated monte carlo assessment factor calculations
$ws.Range("B2").Value = 0.00003519930706972955
$ws.Range("G2").Value = 3.468721592776412
$ws.Range("K2").Value = 0.0196
